$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / first worksheet)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 4380
$ws1.Range("F3").Value = 2458
$ws1.Range("F6").Value = 52
$ws1.Range("F8").Value = 220
$ws1.Range("F9").Value = 132
$ws1.Range("F10").Value = 152
$ws1.Range("F11").Value = 162
$ws1.Range("F12").Value = 1635
$ws1.Range("F13").Value = 298
$ws1.Range("F14").Value = 3472

# Sheet "全部类型" (sheet4 / fourth worksheet)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 4380
$ws4.Range("F3").Value = 2458
$ws4.Range("F7").Value = 52
$ws4.Range("F10").Value = 220
$ws4.Range("F11").Value = 132
$ws4.Range("F12").Value = 152
$ws4.Range("F13").Value = 162
$ws4.Range("F16").Value = 1635
$ws4.Range("F17").Value = 298
$ws4.Range("F18").Value = 3472
